$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 589 (shifts existing rows 589..690 down to 591..692)
$ws.Rows("589:590").Insert()

# --- New row 589 ---
$ws.Cells.Item(589, 1).Value2  = 10
$ws.Cells.Item(589, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(589, 3).Value2  = "La Araucanía"
$ws.Cells.Item(589, 4).Value2  = 44951
$ws.Cells.Item(589, 5).Value2  = 9
$ws.Cells.Item(589, 6).Value2  = 100112043
$ws.Cells.Item(589, 7).Value2  = "Pepino ensalada"
$ws.Cells.Item(589, 8).Value2  = "Sin especificar"
$ws.Cells.Item(589, 9).Value2  = "Primera"
$ws.Cells.Item(589, 10).Value2 = 125
$ws.Cells.Item(589, 11).Value2 = 10000
$ws.Cells.Item(589, 12).Value2 = 10000
$ws.Cells.Item(589, 13).Value2 = 10000
$ws.Cells.Item(589, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(589, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(589, 16).Value2 = 167
$ws.Cells.Item(589, 17).Value2 = 60
$ws.Cells.Item(589, 18).Value2 = "Hortaliza"

# --- New row 590 ---
$ws.Cells.Item(590, 1).Value2  = 10
$ws.Cells.Item(590, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(590, 3).Value2  = "La Araucanía"
$ws.Cells.Item(590, 4).Value2  = 44951
$ws.Cells.Item(590, 5).Value2  = 9
$ws.Cells.Item(590, 6).Value2  = 100112043
$ws.Cells.Item(590, 7).Value2  = "Pepino ensalada"
$ws.Cells.Item(590, 8).Value2  = "Sin especificar"
$ws.Cells.Item(590, 9).Value2  = "Primera"
$ws.Cells.Item(590, 10).Value2 = 280
$ws.Cells.Item(590, 11).Value2 = 12000
$ws.Cells.Item(590, 12).Value2 = 13000
$ws.Cells.Item(590, 13).Value2 = 12304
$ws.Cells.Item(590, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(590, 15).Value2 = "Región del Maule"
$ws.Cells.Item(590, 16).Value2 = 205
$ws.Cells.Item(590, 17).Value2 = 60
$ws.Cells.Item(590, 18).Value2 = "Hortaliza"
